# Remove the "Lifestyle factors" section (header row + Overweight/Obesity/
# Smoking rows) from the end of the MACE regression table, per the commit:
# "Created code for regression of cardiovascular risk factors and
#  eliminated lifestyle factors from MACE survival analysis"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Walk backwards from the last row, deleting the trailing "Lifestyle
# factors" block (the section header cell plus its three data rows) so
# the operation is robust to exactly how many rows precede it.
$cutIndex = -1
for ($i = $t.Rows.Count; $i -ge 1; $i--) {
    $label = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($label -match "Lifestyle factors") {
        $cutIndex = $i
        break
    }
}

if ($cutIndex -gt 0) {
    for ($i = $t.Rows.Count; $i -ge $cutIndex; $i--) {
        $t.Rows.Item($i).Delete()
    }
}
